$d = $word.ActiveDocument

$replacements = @(
    @("2024-12-30 Monday", "2024-12-31 Tuesday"),
    @("13×12=", "12×89="),
    @("39×67=", "66×55="),
    @("98×46=", "50×54="),
    @("62×79=", "58×99="),
    @("17×48=", "11×88="),
    @("49×21=", "36×84="),
    @("20×51=", "52×11="),
    @("59×49=", "93×19="),
    @("36×64=", "86×23="),
    @("96×89=", "94×13="),
    @("77×29=", "58×71="),
    @("68×86=", "72×66="),
    @("99×89=", "59×74="),
    @("94×22=", "80×26="),
    @("32×93=", "73×84="),
    @("88×16=", "80×57="),
    @("23×82=", "49×54="),
    @("26×49=", "43×96="),
    @("23×95=", "90×53="),
    @("19×72=", "60×54="),
    @("34×54=", "17×73="),
    @("52×80=", "61×56="),
    @("36×23=", "81×38="),
    @("14×24=", "58×60="),
    @("59×82=", "35×84=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
